# Add new rows describing additional cardio-metabolic outcome variables
# (IDEFICS reference-based z-scores / percentiles) to the "Variables" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

$newRows = @(
    @("bmi_idefics_perc_", "decimal", "Repeated measures of child's BMI z-score percentiles based on the IDEFICS study reference values"),
    @("bmi_idefics_z_", "decimal", "Repeated measures of child's BMI z-score  based on the IDEFICS study reference values"),
    @("wc_idefics_z_", "decimal", "Repeated measures of child's waist circumference z-score  based on the IDEFICS study reference values"),
    @("hdl_idefics_z_", "decimal", "Repeated measures of child's HDL z-score  based on the IDEFICS study reference values"),
    @("triglyceride_idefics_z_", "decimal", "Repeated measures of child's triglycerides z-score  based on the IDEFICS study reference values"),
    @("glucose_idefics_z_", "decimal", "Repeated measures of child's glucose z-score  based on the IDEFICS study reference values"),
    @("sbp_idefics_z_", "decimal", "Repeated measures of child's SBP z-score  based on the IDEFICS study reference values"),
    @("dbp_idefics_z_", "decimal", "Repeated measures of child's DBP z-score  based on the IDEFICS study reference values"),
    @("mets_nriskfactors_idefics_", "integer", "risk factors based on the IDEFICS study reference values")
)

$startRow = 88
$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}

$lastRow = $r - 1
$ws.Range("A$($startRow):D$($lastRow)").Select()
